# Fill in columns G:N with "N/A" for the rows that were previously
# missing those cells (rows that only had data through column F).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(5,6,11,13,15,17,20,23,30,34,38,40,41,45,46,48,52,53,54,58,59,63,64,65,71,72,76,78,79,85,90,92,97,101,107,110,112,113,114,116,118,119,122,124,127,128,133,139,144,145,151,152,156,160,161,164,166,171,175,177,179,184,190,192)

foreach ($r in $rows) {
    $ws.Range("G" + $r + ":N" + $r).Value = "N/A"
}
